$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert 4 new blank rows before row 178 (each single-row insert shifts
# everything below down by one and inherits formatting from the row above,
# matching the A/E/G column styling already used through this block).
$ws.Rows.Item(178).Insert(-4121)
$ws.Rows.Item(178).Insert(-4121)
$ws.Rows.Item(178).Insert(-4121)
$ws.Rows.Item(178).Insert(-4121)

# Grow the worksheet Table (Table1) to cover the 4 extra rows.
$lo.Resize($ws.Range("A1:G274"))

# Fill in A/B ("set"/"subset") for the 4 new rows first - these reuse
# already-existing shared strings ("@C#" / "Basics").
$ws.Range("A178").Value = "'@C#"
$ws.Range("B178").Value = "Basics"
$ws.Range("A179").Value = "'@C#"
$ws.Range("B179").Value = "Basics"
$ws.Range("A180").Value = "'@C#"
$ws.Range("B180").Value = "Basics"
$ws.Range("A181").Value = "'@C#"
$ws.Range("B181").Value = "Basics"

# Column F ("links") filled top to bottom.
$ws.Range("F178").Value = "Render separate cshtml (sub-template)"
$ws.Range("F179").Value = "Render separate cshtml with params (sub-template)"
$ws.Range("F180").Value = "Create instance of CSHTML with library functions"
$ws.Range("F181").Value = "In a sub-template, access parameter given in"

# Column C ("name") - rows 180/181 reuse the text already used in F180/F181.
$ws.Range("C178").Value = "Render sub-template"
$ws.Range("C179").Value = "Render sub-template w/params"
$ws.Range("C180").Value = "Create instance of CSHTML with library functions"
$ws.Range("C181").Value = "In a sub-template, access parameter given in"

# Column E ("content") - the actual snippet bodies, entered 180, 181, 179, 178.
$e180 = "'@{" + "`n" + '    var lib = CreateInstance("_${1:library}.cshtml"); ' + "`n" + '}'
$ws.Range("E180").Value = $e180

$e181 = "'@{" + "`n" + '    var ${2:post} = PageData["${1:Post}"];' + "`n" + '}'
$ws.Range("E181").Value = $e181

$ws.Range("E179").Value = "'@RenderPage(`"_`${1:list-item}.cshtml`", new { `${2:Post} = `${3:post} })"

$ws.Range("E178").Value = "'@RenderPage(`"_`${1:list-item}.cshtml`")"

# Row heights for the wrapped-text rows.
$ws.Rows.Item(179).RowHeight = 30
$ws.Rows.Item(180).RowHeight = 45
$ws.Rows.Item(181).RowHeight = 45

# Reflect where the author ended up scrolled to / selecting after the edit.
$ws.Application.Goto($ws.Range("A174"))
$ws.Range("E179").Select()
